$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Full Timetable")

$ws.Range("B1").Value = "14 July (Monday)"
$ws.Range("C1").Value = "15 July (Tuesday)"
$ws.Range("D1").Value = "16 July (Wednesday)"
$ws.Range("E1").Value = "17 July (Thursday)"
$ws.Range("F1").Value = "18 July (Friday)"
$ws.Range("G1").Value = "19 July (Saturday)"
